# "Anpassung return + Vortrag"
# Slide 4 ("Modell") contains a BibTeX-field reference table. The last
# column lists field names for the @inproceedings / @book-ish entry type
# ("title", "booktitle", "series", "publisher"). The author extended the
# list with three more BibTeX field names ("Publisher", "Journal",
# "school") and renamed the former "publisher" entry to "month".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(5)          # "Tabelle 14" graphicFrame
$tbl = $shp.Table
$cell = $tbl.Cell(1, 4)           # column holding title/booktitle/series/publisher
$tr = $cell.Shape.TextFrame.TextRange

$tr.Text = "title`rbooktitle`rseries`rPublisher`rJournal`rschool`rmonth"
